# household_member.xlsx update:
#  - update handlebars-style prompts on the "survey" sheet to reference {{member_name}}
#  - add a new "calculates" sheet defining the ageIsOddOrEven calculation
#  - insert a new "note" row that shows a calculate-backed message
#  - tweak selections on a couple of sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. survey sheet: reword prompts to be parameterized by {{member_name}}
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

$survey.Range("D4").Value = "Enter age of {{member_name}}:"
$survey.Range("D5").Value = "Enter sex of {{member_name}}:"
$survey.Range("D6").Value = "Does {{member_name}} contribute to the household income?"

# ---------------------------------------------------------------------------
# 2. new "calculates" sheet, appended after "model" -- add headers first
# ---------------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")
$calculates = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $model)
$calculates.Name = "calculates"

$calculates.Cells.Item(1, 1).Value = "calculation_name"
$calculates.Cells.Item(1, 2).Value = "calculation"

# ---------------------------------------------------------------------------
# 3. back to survey: insert a "note" row above income_contribution, shifting
#    income_contribution (and everything below it) down by one row. The
#    note's message text is filled in later, once the calculation exists.
# ---------------------------------------------------------------------------
$survey.Rows.Item(6).Insert()
$survey.Rows.Item(6).RowHeight = 31
$survey.Cells.Item(6, 1).Value = "note"

# ---------------------------------------------------------------------------
# 4. finish the "calculates" sheet definition
# ---------------------------------------------------------------------------
$calculates.Cells.Item(2, 1).Value = "ageIsOddOrEven"
$calculates.Cells.Item(2, 2).Value = "((data('age') % 2) == 1) ? ""odd"" : ""even"""

$calculates.Range("A2").Select()

# ---------------------------------------------------------------------------
# 5. fill in the note's message, now that the calculation it references
#    exists.
# ---------------------------------------------------------------------------
$survey.Cells.Item(6, 4).Value = "{{member_name}} age is {{evaluate calculates.ageIsOddOrEven}} in {{setting 'table_id'}} for {{metadata 'instanceName'}}"

# ---------------------------------------------------------------------------
# 6. choices sheet: selection moved to the (previously unused) E column
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")
$choices.Columns.Item(5).Select()

# ---------------------------------------------------------------------------
# 7. leave "survey" as the active tab, selection on the income_contribution
#    question (now on row 7).
# ---------------------------------------------------------------------------
$survey.Select()
$survey.Range("D7").Select()
